$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.722.80'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +2.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.851.31'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +1.74%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9978'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.70'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6384'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +4.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9984'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.56%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.846.72'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +1.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07498'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +2.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2966'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +3.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.93'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +5.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07676'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +0.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.833.63'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.063'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +3.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6881'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +5.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.81'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +4.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009375'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +5.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.027'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +3.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.660.69'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +2.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.084.41'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +0.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '240.58'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +2.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.65'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +2.37%  '

$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.392'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +4.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.9992'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.22'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +0.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1429'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +3.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.560'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +2.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.94'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +2.28%  '

$ws.Range("E30").Value = '  +0.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.06045'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +8.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.259'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +4.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.160'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +2.68%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.146'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +2.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.881'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +3.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.153'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +2.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7325'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +1.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.606'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -0.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.877'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +2.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.233.06'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +3.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01777'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +1.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.350'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +0.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9179'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +4.15%  '

$ws.Range("E44").Value = '  -0.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.003.09'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.10'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +1.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.45'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +4.04%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000120'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5078'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.297'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +3.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4095'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +3.33%  '
